$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 28.56536566666667
$ws.Cells.Item(2, 8).Value = 85.69609700000001
$ws.Cells.Item(2, 9).Value = 0.02097368575335975
$ws.Cells.Item(2, 10).Value = 0.02097368575335974
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 0.043534
$ws.Cells.Item(2, 14).Value = 0.130602
$ws.Cells.Item(2, 15).Value = 0.00760428821720775
$ws.Cells.Item(2, 16).Value = 0.007604288217207752
$ws.Cells.Item(2, 17).Value = 1.243564628932667
$ws.Cells.Item(2, 18).Value = 11.192081660394
$ws.Cells.Item(2, 19).Value = 0.0001594899514456916
$ws.Cells.Item(2, 20).Value = 0.0001594899514456916

$ws.Cells.Item(3, 7).Value = 28.56536566666667
$ws.Cells.Item(3, 8).Value = 85.69609700000001
$ws.Cells.Item(3, 9).Value = 0.02097368575335975
$ws.Cells.Item(3, 10).Value = 0.02097368575335974
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 5.524447333333334
$ws.Cells.Item(3, 14).Value = 16.573342
$ws.Cells.Item(3, 15).Value = 0.964981158713912
$ws.Cells.Item(3, 16).Value = 0.9649811587139122
$ws.Cells.Item(3, 17).Value = 157.8078581829082
$ws.Cells.Item(3, 18).Value = 1420.270723646174
$ws.Cells.Item(3, 19).Value = 0.02023921158077856
$ws.Cells.Item(3, 20).Value = 0.02023921158077856

$ws.Cells.Item(4, 7).Value = 28.56536566666667
$ws.Cells.Item(4, 8).Value = 85.69609700000001
$ws.Cells.Item(4, 9).Value = 0.02097368575335975
$ws.Cells.Item(4, 10).Value = 0.02097368575335974
$ws.Cells.Item(4, 11).Value = 1
$ws.Cells.Item(4, 12).Value = 0.3333333333333333
$ws.Cells.Item(4, 13).Value = 0.03706533333333333
$ws.Cells.Item(4, 14).Value = 0.111196
$ws.Cells.Item(4, 15).Value = 0.006474375833453032
$ws.Cells.Item(4, 16).Value = 0.006474375833453035
$ws.Cells.Item(4, 17).Value = 1.058784800223556
$ws.Cells.Item(4, 18).Value = 9.529063202012
$ws.Cells.Item(4, 19).Value = 0.0001357915241799905
$ws.Cells.Item(4, 20).Value = 0.0001357915241799906

$ws.Cells.Item(5, 7).Value = 28.56536566666667
$ws.Cells.Item(5, 8).Value = 85.69609700000001
$ws.Cells.Item(5, 9).Value = 0.02097368575335975
$ws.Cells.Item(5, 10).Value = 0.02097368575335974
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 0.119881
$ws.Cells.Item(5, 14).Value = 0.359643
$ws.Cells.Item(5, 15).Value = 0.02094017723542708
$ws.Cells.Item(5, 16).Value = 0.02094017723542708
$ws.Cells.Item(5, 17).Value = 3.424444601485667
$ws.Cells.Item(5, 18).Value = 30.82000141337101
$ws.Cells.Item(5, 19).Value = 0.000439192696955505
$ws.Cells.Item(5, 20).Value = 0.0004391926969555051

$ws.Cells.Item(6, 7).Value = 1288.520629666667
$ws.Cells.Item(6, 8).Value = 3865.561889000001
$ws.Cells.Item(6, 9).Value = 0.9460766961189575
$ws.Cells.Item(6, 10).Value = 0.9460766961189573
$ws.Cells.Item(6, 11).Value = 1
$ws.Cells.Item(6, 12).Value = 0.3333333333333333
$ws.Cells.Item(6, 13).Value = 0.043534
$ws.Cells.Item(6, 14).Value = 0.130602
$ws.Cells.Item(6, 15).Value = 0.00760428821720775
$ws.Cells.Item(6, 16).Value = 0.007604288217207752
$ws.Cells.Item(6, 17).Value = 56.09445709190867
$ws.Cells.Item(6, 18).Value = 504.8501138271781
$ws.Cells.Item(6, 19).Value = 0.007194239872872226
$ws.Cells.Item(6, 20).Value = 0.007194239872872227

$ws.Cells.Item(7, 7).Value = 1288.520629666667
$ws.Cells.Item(7, 8).Value = 3865.561889000001
$ws.Cells.Item(7, 9).Value = 0.9460766961189575
$ws.Cells.Item(7, 10).Value = 0.9460766961189573
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 5.524447333333334
$ws.Cells.Item(7, 14).Value = 16.573342
$ws.Cells.Item(7, 15).Value = 0.964981158713912
$ws.Cells.Item(7, 16).Value = 0.9649811587139122
$ws.Cells.Item(7, 17).Value = 7118.364356507005
$ws.Cells.Item(7, 18).Value = 64065.27920856305
$ws.Cells.Item(7, 19).Value = 0.9129461864531012
$ws.Cells.Item(7, 20).Value = 0.9129461864531012

$ws.Cells.Item(8, 7).Value = 1288.520629666667
$ws.Cells.Item(8, 8).Value = 3865.561889000001
$ws.Cells.Item(8, 9).Value = 0.9460766961189575
$ws.Cells.Item(8, 10).Value = 0.9460766961189573
$ws.Cells.Item(8, 11).Value = 1
$ws.Cells.Item(8, 12).Value = 0.3333333333333333
$ws.Cells.Item(8, 13).Value = 0.03706533333333333
$ws.Cells.Item(8, 14).Value = 0.111196
$ws.Cells.Item(8, 15).Value = 0.006474375833453032
$ws.Cells.Item(8, 16).Value = 0.006474375833453035
$ws.Cells.Item(8, 17).Value = 47.75944664547156
$ws.Cells.Item(8, 18).Value = 429.835019809244
$ws.Cells.Item(8, 19).Value = 0.006125256097945667
$ws.Cells.Item(8, 20).Value = 0.006125256097945668

$ws.Cells.Item(9, 7).Value = 1288.520629666667
$ws.Cells.Item(9, 8).Value = 3865.561889000001
$ws.Cells.Item(9, 9).Value = 0.9460766961189575
$ws.Cells.Item(9, 10).Value = 0.9460766961189573
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 0.119881
$ws.Cells.Item(9, 14).Value = 0.359643
$ws.Cells.Item(9, 15).Value = 0.02094017723542708
$ws.Cells.Item(9, 16).Value = 0.02094017723542708
$ws.Cells.Item(9, 17).Value = 154.4691416050697
$ws.Cells.Item(9, 18).Value = 1390.222274445627
$ws.Cells.Item(9, 19).Value = 0.01981101369503825
$ws.Cells.Item(9, 20).Value = 0.01981101369503826

$ws.Cells.Item(10, 7).Value = 0.115045
$ws.Cells.Item(10, 8).Value = 0.345135
$ws.Cells.Item(10, 9).Value = 0.00008447004339632664
$ws.Cells.Item(10, 10).Value = 0.00008447004339632662
$ws.Cells.Item(10, 11).Value = 1
$ws.Cells.Item(10, 12).Value = 0.3333333333333333
$ws.Cells.Item(10, 13).Value = 0.043534
$ws.Cells.Item(10, 14).Value = 0.130602
$ws.Cells.Item(10, 15).Value = 0.00760428821720775
$ws.Cells.Item(10, 16).Value = 0.007604288217207752
$ws.Cells.Item(10, 17).Value = 0.00500836903
$ws.Cells.Item(10, 18).Value = 0.04507532127
$ws.Cells.Item(10, 19).Value = 0.0000006423345557057139
$ws.Cells.Item(10, 20).Value = 0.000000642334555705714

$ws.Cells.Item(11, 7).Value = 0.115045
$ws.Cells.Item(11, 8).Value = 0.345135
$ws.Cells.Item(11, 9).Value = 0.00008447004339632664
$ws.Cells.Item(11, 10).Value = 0.00008447004339632662
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 5.524447333333334
$ws.Cells.Item(11, 14).Value = 16.573342
$ws.Cells.Item(11, 15).Value = 0.964981158713912
$ws.Cells.Item(11, 16).Value = 0.9649811587139122
$ws.Cells.Item(11, 17).Value = 0.6355600434633334
$ws.Cells.Item(11, 18).Value = 5.72004039117
$ws.Cells.Item(11, 19).Value = 0.00008151200035320171
$ws.Cells.Item(11, 20).Value = 0.00008151200035320171

$ws.Cells.Item(12, 7).Value = 0.115045
$ws.Cells.Item(12, 8).Value = 0.345135
$ws.Cells.Item(12, 9).Value = 0.00008447004339632664
$ws.Cells.Item(12, 10).Value = 0.00008447004339632662
$ws.Cells.Item(12, 11).Value = 1
$ws.Cells.Item(12, 12).Value = 0.3333333333333333
$ws.Cells.Item(12, 13).Value = 0.03706533333333333
$ws.Cells.Item(12, 14).Value = 0.111196
$ws.Cells.Item(12, 15).Value = 0.006474375833453032
$ws.Cells.Item(12, 16).Value = 0.006474375833453035
$ws.Cells.Item(12, 17).Value = 0.004264181273333333
$ws.Cells.Item(12, 18).Value = 0.03837763146000001
$ws.Cells.Item(12, 19).Value = 0.000000546890807615906
$ws.Cells.Item(12, 20).Value = 0.0000005468908076159063

$ws.Cells.Item(13, 7).Value = 0.115045
$ws.Cells.Item(13, 8).Value = 0.345135
$ws.Cells.Item(13, 9).Value = 0.00008447004339632664
$ws.Cells.Item(13, 10).Value = 0.00008447004339632662
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 0.119881
$ws.Cells.Item(13, 14).Value = 0.359643
$ws.Cells.Item(13, 15).Value = 0.02094017723542708
$ws.Cells.Item(13, 16).Value = 0.02094017723542708
$ws.Cells.Item(13, 17).Value = 0.013791709645
$ws.Cells.Item(13, 18).Value = 0.124125386805
$ws.Cells.Item(13, 19).Value = 0.000001768817679803296
$ws.Cells.Item(13, 20).Value = 0.000001768817679803297

$ws.Cells.Item(14, 7).Value = 43.90798866666668
$ws.Cells.Item(14, 8).Value = 131.723966
$ws.Cells.Item(14, 9).Value = 0.03223877359397412
$ws.Cells.Item(14, 10).Value = 0.0322387735939741
$ws.Cells.Item(14, 11).Value = 1
$ws.Cells.Item(14, 12).Value = 0.3333333333333333
$ws.Cells.Item(14, 13).Value = 0.043534
$ws.Cells.Item(14, 14).Value = 0.130602
$ws.Cells.Item(14, 15).Value = 0.00760428821720775
$ws.Cells.Item(14, 16).Value = 0.007604288217207752
$ws.Cells.Item(14, 17).Value = 1.911490378614667
$ws.Cells.Item(14, 18).Value = 17.203413407532
$ws.Cells.Item(14, 19).Value = 0.0002451529261778857
$ws.Cells.Item(14, 20).Value = 0.0002451529261778857

$ws.Cells.Item(15, 7).Value = 43.90798866666668
$ws.Cells.Item(15, 8).Value = 131.723966
$ws.Cells.Item(15, 9).Value = 0.03223877359397412
$ws.Cells.Item(15, 10).Value = 0.0322387735939741
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 5.524447333333334
$ws.Cells.Item(15, 14).Value = 16.573342
$ws.Cells.Item(15, 15).Value = 0.964981158713912
$ws.Cells.Item(15, 16).Value = 0.9649811587139122
$ws.Cells.Item(15, 17).Value = 242.5673709015969
$ws.Cells.Item(15, 18).Value = 2183.106338114372
$ws.Cells.Item(15, 19).Value = 0.03110980909822861
$ws.Cells.Item(15, 20).Value = 0.03110980909822861

$ws.Cells.Item(16, 7).Value = 43.90798866666668
$ws.Cells.Item(16, 8).Value = 131.723966
$ws.Cells.Item(16, 9).Value = 0.03223877359397412
$ws.Cells.Item(16, 10).Value = 0.0322387735939741
$ws.Cells.Item(16, 11).Value = 1
$ws.Cells.Item(16, 12).Value = 0.3333333333333333
$ws.Cells.Item(16, 13).Value = 0.03706533333333333
$ws.Cells.Item(16, 14).Value = 0.111196
$ws.Cells.Item(16, 15).Value = 0.006474375833453032
$ws.Cells.Item(16, 16).Value = 0.006474375833453035
$ws.Cells.Item(16, 17).Value = 1.627464235926223
$ws.Cells.Item(16, 18).Value = 14.647178123336
$ws.Cells.Item(16, 19).Value = 0.0002087259366569898
$ws.Cells.Item(16, 20).Value = 0.0002087259366569898

$ws.Cells.Item(17, 7).Value = 43.90798866666668
$ws.Cells.Item(17, 8).Value = 131.723966
$ws.Cells.Item(17, 9).Value = 0.03223877359397412
$ws.Cells.Item(17, 10).Value = 0.0322387735939741
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 0.119881
$ws.Cells.Item(17, 14).Value = 0.359643
$ws.Cells.Item(17, 15).Value = 0.02094017723542708
$ws.Cells.Item(17, 16).Value = 0.02094017723542708
$ws.Cells.Item(17, 17).Value = 5.263733589348669
$ws.Cells.Item(17, 18).Value = 47.37360230413801
$ws.Cells.Item(17, 19).Value = 0.0006750856329106245
$ws.Cells.Item(17, 20).Value = 0.0006750856329106244

$ws.Cells.Item(18, 7).Value = 0.693788
$ws.Cells.Item(18, 8).Value = 2.081364
$ws.Cells.Item(18, 9).Value = 0.0005094032984297506
$ws.Cells.Item(18, 10).Value = 0.0005094032984297505
$ws.Cells.Item(18, 11).Value = 1
$ws.Cells.Item(18, 12).Value = 0.3333333333333333
$ws.Cells.Item(18, 13).Value = 0.043534
$ws.Cells.Item(18, 14).Value = 0.130602
$ws.Cells.Item(18, 15).Value = 0.00760428821720775
$ws.Cells.Item(18, 16).Value = 0.007604288217207752
$ws.Cells.Item(18, 17).Value = 0.030203366792
$ws.Cells.Item(18, 18).Value = 0.2718303011279999
$ws.Cells.Item(18, 19).Value = 0.000003873649500056116
$ws.Cells.Item(18, 20).Value = 0.000003873649500056116

$ws.Cells.Item(19, 7).Value = 0.693788
$ws.Cells.Item(19, 8).Value = 2.081364
$ws.Cells.Item(19, 9).Value = 0.0005094032984297506
$ws.Cells.Item(19, 10).Value = 0.0005094032984297505
$ws.Cells.Item(19, 11).Value = 3
$ws.Cells.Item(19, 12).Value = 1
$ws.Cells.Item(19, 13).Value = 5.524447333333334
$ws.Cells.Item(19, 14).Value = 16.573342
$ws.Cells.Item(19, 15).Value = 0.964981158713912
$ws.Cells.Item(19, 16).Value = 0.9649811587139122
$ws.Cells.Item(19, 17).Value = 3.832795266498667
$ws.Cells.Item(19, 18).Value = 34.495157398488
$ws.Cells.Item(19, 19).Value = 0.0004915645851714294
$ws.Cells.Item(19, 20).Value = 0.0004915645851714294

$ws.Cells.Item(20, 7).Value = 0.693788
$ws.Cells.Item(20, 8).Value = 2.081364
$ws.Cells.Item(20, 9).Value = 0.0005094032984297506
$ws.Cells.Item(20, 10).Value = 0.0005094032984297505
$ws.Cells.Item(20, 11).Value = 1
$ws.Cells.Item(20, 12).Value = 0.3333333333333333
$ws.Cells.Item(20, 13).Value = 0.03706533333333333
$ws.Cells.Item(20, 14).Value = 0.111196
$ws.Cells.Item(20, 15).Value = 0.006474375833453032
$ws.Cells.Item(20, 16).Value = 0.006474375833453035
$ws.Cells.Item(20, 17).Value = 0.02571548348266666
$ws.Cells.Item(20, 18).Value = 0.231439351344
$ws.Cells.Item(20, 19).Value = 0.00000329806840483484
$ws.Cells.Item(20, 20).Value = 0.000003298068404834841

$ws.Cells.Item(21, 7).Value = 0.693788
$ws.Cells.Item(21, 8).Value = 2.081364
$ws.Cells.Item(21, 9).Value = 0.0005094032984297506
$ws.Cells.Item(21, 10).Value = 0.0005094032984297505
$ws.Cells.Item(21, 11).Value = 3
$ws.Cells.Item(21, 12).Value = 1
$ws.Cells.Item(21, 13).Value = 0.119881
$ws.Cells.Item(21, 14).Value = 0.359643
$ws.Cells.Item(21, 15).Value = 0.02094017723542708
$ws.Cells.Item(21, 16).Value = 0.02094017723542708
$ws.Cells.Item(21, 17).Value = 0.083171999228
$ws.Cells.Item(21, 18).Value = 0.748547993052
$ws.Cells.Item(21, 19).Value = 0.00001066699535343013
$ws.Cells.Item(21, 20).Value = 0.00001066699535343013

$ws.Cells.Item(22, 7).Value = 0.1593103333333333
$ws.Cells.Item(22, 8).Value = 0.477931
$ws.Cells.Item(22, 9).Value = 0.0001169711918827409
$ws.Cells.Item(22, 10).Value = 0.0001169711918827409
$ws.Cells.Item(22, 11).Value = 1
$ws.Cells.Item(22, 12).Value = 0.3333333333333333
$ws.Cells.Item(22, 13).Value = 0.043534
$ws.Cells.Item(22, 14).Value = 0.130602
$ws.Cells.Item(22, 15).Value = 0.00760428821720775
$ws.Cells.Item(22, 16).Value = 0.007604288217207752
$ws.Cells.Item(22, 17).Value = 0.006935416051333332
$ws.Cells.Item(22, 18).Value = 0.062418744462
$ws.Cells.Item(22, 19).Value = 0.0000008894826561866736
$ws.Cells.Item(22, 20).Value = 0.0000008894826561866737

$ws.Cells.Item(23, 7).Value = 0.1593103333333333
$ws.Cells.Item(23, 8).Value = 0.477931
$ws.Cells.Item(23, 9).Value = 0.0001169711918827409
$ws.Cells.Item(23, 10).Value = 0.0001169711918827409
$ws.Cells.Item(23, 11).Value = 3
$ws.Cells.Item(23, 12).Value = 1
$ws.Cells.Item(23, 13).Value = 5.524447333333334
$ws.Cells.Item(23, 14).Value = 16.573342
$ws.Cells.Item(23, 15).Value = 0.964981158713912
$ws.Cells.Item(23, 16).Value = 0.9649811587139122
$ws.Cells.Item(23, 17).Value = 0.8801015461557778
$ws.Cells.Item(23, 18).Value = 7.920913915402
$ws.Cells.Item(23, 19).Value = 0.0001128749962791547
$ws.Cells.Item(23, 20).Value = 0.0001128749962791547

$ws.Cells.Item(24, 7).Value = 0.1593103333333333
$ws.Cells.Item(24, 8).Value = 0.477931
$ws.Cells.Item(24, 9).Value = 0.0001169711918827409
$ws.Cells.Item(24, 10).Value = 0.0001169711918827409
$ws.Cells.Item(24, 11).Value = 1
$ws.Cells.Item(24, 12).Value = 0.3333333333333333
$ws.Cells.Item(24, 13).Value = 0.03706533333333333
$ws.Cells.Item(24, 14).Value = 0.111196
$ws.Cells.Item(24, 15).Value = 0.006474375833453032
$ws.Cells.Item(24, 16).Value = 0.006474375833453035
$ws.Cells.Item(24, 17).Value = 0.005904890608444444
$ws.Cells.Item(24, 18).Value = 0.053144015476
$ws.Cells.Item(24, 19).Value = 0.0000007573154579358152
$ws.Cells.Item(24, 20).Value = 0.0000007573154579358154

$ws.Cells.Item(25, 7).Value = 0.1593103333333333
$ws.Cells.Item(25, 8).Value = 0.477931
$ws.Cells.Item(25, 9).Value = 0.0001169711918827409
$ws.Cells.Item(25, 10).Value = 0.0001169711918827409
$ws.Cells.Item(25, 11).Value = 3
$ws.Cells.Item(25, 12).Value = 1
$ws.Cells.Item(25, 13).Value = 0.119881
$ws.Cells.Item(25, 14).Value = 0.359643
$ws.Cells.Item(25, 15).Value = 0.02094017723542708
$ws.Cells.Item(25, 16).Value = 0.02094017723542708
$ws.Cells.Item(25, 17).Value = 0.01909828207033333
$ws.Cells.Item(25, 18).Value = 0.171884538633
$ws.Cells.Item(25, 19).Value = 0.000002449397489463744
$ws.Cells.Item(25, 20).Value = 0.000002449397489463744

Write-Output "Updated Col1a1-Itga11 TPM values"